# "time is now static" / "minor bug fix" - the source data-generation run
# re-exported this factory file: the duplicated, leftover "charcoal" power
# row on the connections sheet is cleared out (only its number-format
# styling is kept), and the saved cursor/selection on a couple of sheets
# moved since the export now stops right after the last real row.

$wb = $excel.ActiveWorkbook

# connections sheet: row 26 was an accidental duplicate of row 25's
# "power" connection (for charcoal) - drop its values, keep styles.
$wsConnections = $wb.Worksheets.Item("connections")
$wsConnections.Range("B26:I26").ClearContents() | Out-Null
$wsConnections.Range("C30").Select() | Out-Null

# chains sheet: selection left on the row below the used range.
$wsChains = $wb.Worksheets.Item("chains")
$wsChains.Range("A12:E12").Select() | Out-Null

# leave "connections" as the active/selected sheet, matching activeTab=1
$wsConnections.Activate() | Out-Null
$wsConnections.Range("C30").Select() | Out-Null
